# Bundle rm_sumrow() and rm_sumcol() to rebel():
# Add a "contamination" sheet (Sheet2) ahead of the existing Sheet1, containing
# the same layout/labels/formulas but with every data value shifted by +210.

$wb = $excel.ActiveWorkbook

# Add the new sheet. Worksheets.Add() inserts before the active sheet and
# becomes the new active sheet itself -- Excel auto-names it "Sheet2" since
# "Sheet1" already exists, giving the exact "Sheet2", "Sheet1" tab order.
$ws2 = $wb.Worksheets.Add()

# NOTE: fetch this *after* Add() -- a reference grabbed before adding a sheet
# goes stale (reads back empty) once the sheet collection has changed.
$ws1 = $wb.Worksheets.Item("Sheet1")

$dataRows = @(2,3,4,5,6,7,8,9,10,12,13,14,15,16,17,18,19,20,22,23,24,25,26,27,28,29,30)
$cols = @("B","C","D","E","F","G","H")

# Header row (A1:H1) -- identical labels to Sheet1.
foreach ($col in @("A","B","C","D","E","F","G","H")) {
    $addr = $col + "1"
    $ws2.Range($addr).Value = $ws1.Range($addr).Value2
}

# Column-A row labels for every data row, and the "sum" labels on rows 11/21.
foreach ($r in $dataRows) {
    $addr = "A" + $r
    $ws2.Range($addr).Value = $ws1.Range($addr).Value2
}
$ws2.Range("A11").Value = $ws1.Range("A11").Value2
$ws2.Range("A21").Value = $ws1.Range("A21").Value2

# Data cells B:H for every data row, shifted by +210 vs. Sheet1.
foreach ($r in $dataRows) {
    foreach ($col in $cols) {
        $addr = $col + $r
        $ws2.Range($addr).Value = $ws1.Range($addr).Value2 + 210
    }
}

# Sum formulas, matching Sheet1's layout (row 11 sums rows 2:10, row 21 sums rows 12:20).
foreach ($col in $cols) {
    $ws2.Range($col + "11").Formula = "=SUM(" + $col + "2:" + $col + "10)"
    $ws2.Range($col + "21").Formula = "=SUM(" + $col + "12:" + $col + "20)"
}

# Cursor/selection left on the new sheet at J18 (an empty cell past the data),
# and that sheet stays the active tab -- matching the saved selection state.
$ws2.Range("J18").Select() | Out-Null
